$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-12 Sunday" "2025-10-13 Monday"

Replace-Text "56÷3=18, 2" "81÷7=11, 4"
Replace-Text "41÷6=6, 5" "15÷2=7, 1"
Replace-Text "49÷5=9, 4" "28÷3=9, 1"
Replace-Text "85÷4=21, 1" "66÷9=7, 3"
Replace-Text "18÷7=2, 4" "53÷5=10, 3"

Replace-Text "25÷8=3, 1" "39÷9=4, 3"
Replace-Text "92÷8=11, 4" "64÷9=7, 1"
Replace-Text "21÷6=3, 3" "77÷6=12, 5"
Replace-Text "12÷3=4, 0" "43÷5=8, 3"
Replace-Text "61÷8=7, 5" "31÷9=3, 4"

Replace-Text "22÷4=5, 2" "85÷8=10, 5"
Replace-Text "93÷4=23, 1" "87÷3=29, 0"
Replace-Text "48÷7=6, 6" "71÷6=11, 5"
Replace-Text "80÷5=16, 0" "13÷9=1, 4"
Replace-Text "54÷8=6, 6" "74÷8=9, 2"

Replace-Text "55÷7=7, 6" "65÷4=16, 1"
Replace-Text "46÷3=15, 1" "30÷3=10, 0"
Replace-Text "54÷7=7, 5" "84÷9=9, 3"
Replace-Text "56÷9=6, 2" "49÷2=24, 1"
Replace-Text "94÷8=11, 6" "79÷4=19, 3"

Replace-Text "67÷7=9, 4" "95÷3=31, 2"
Replace-Text "34÷3=11, 1" "48÷2=24, 0"
Replace-Text "96÷7=13, 5" "87÷3=29, 0"
Replace-Text "16÷6=2, 4" "69÷2=34, 1"
Replace-Text "90÷5=18, 0" "91÷5=18, 1"
